$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns J and K (Round_7, Round_8) entirely
$ws.Range("J1:K6").Delete() | Out-Null

# Update the remaining score values (D2:I6) to their new values
$values = @(
    @(9,7,2,6,6,7),
    @(7,8,2,1,3,4),
    @(7,5,9,8,3,9),
    @(1,2,4,0,3,5),
    @(0,7,5,2,6,1)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = $j + 4  # D = 4
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}
